# Trade #162 closed at 2026-02-18 00:44:27 - unknown UNKNOWN +0.000%
#
# - Closes the open "momentum" trade (Trade #190) as an early exit.
# - Opens two brand-new trades: #219 (HighProbConvergence) and #220 (MarketMaking).
# - Refreshes the roll-up counters on the Summary / Strategy Status sheets.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 190     # Total Trades
$summary.Range("B9").Value = 43.68   # Win Rate %

# ---------------------------------------------------------------------------
# Strategy Status sheet (row 11 = "momentum" strategy)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D11").Value = 50   # Trades
$status.Range("G11").Value = 30   # Win Rate %

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

# Close out Trade #190 (row 191) - momentum strategy, early exit.
$allTrades.Cells.Item(191, 7).Value = 0.01          # Exit Price
$allTrades.Cells.Item(191, 8).Value = "CLOSED"      # Status
$allTrades.Cells.Item(191, 11).Value = 99.22        # Capital After
$allTrades.Cells.Item(191, 12).Value = "early_exit" # Exit Reason
$allTrades.Cells.Item(191, 13).Value = 0.18         # Duration (min)

# New trade #219 - HighProbConvergence (row 220)
$allTrades.Range("B220").NumberFormat = "@"
$allTrades.Cells.Item(220, 1).Value = 219
$allTrades.Cells.Item(220, 2).Value = "2026-02-18"
$allTrades.Cells.Item(220, 3).Value = "00:44:21"
$allTrades.Cells.Item(220, 4).Value = "HighProbConvergence"
$allTrades.Cells.Item(220, 5).Value = "UP"
$allTrades.Cells.Item(220, 6).Value = 0.01
$allTrades.Cells.Item(220, 8).Value = "OPEN"
$allTrades.Cells.Item(220, 9).Value = 0
$allTrades.Cells.Item(220, 10).Value = 0
$allTrades.Cells.Item(220, 11).Value = 100.3223499536821
$allTrades.Cells.Item(220, 13).Value = 0
$allTrades.Cells.Item(220, 14).Value = 0
$allTrades.Cells.Item(220, 15).Value = 0
$allTrades.Cells.Item(220, 16).Value = 0.95
$allTrades.Cells.Item(220, 17).Value = "Mean reversion UP: price 10.82% below mean (z=-3.00)"

# New trade #220 - MarketMaking (row 221)
$allTrades.Range("B221").NumberFormat = "@"
$allTrades.Cells.Item(221, 1).Value = 220
$allTrades.Cells.Item(221, 2).Value = "2026-02-18"
$allTrades.Cells.Item(221, 3).Value = "00:44:22"
$allTrades.Cells.Item(221, 4).Value = "MarketMaking"
$allTrades.Cells.Item(221, 5).Value = "UP"
$allTrades.Cells.Item(221, 6).Value = 0.01
$allTrades.Cells.Item(221, 8).Value = "OPEN"
$allTrades.Cells.Item(221, 9).Value = 0
$allTrades.Cells.Item(221, 10).Value = 0
$allTrades.Cells.Item(221, 11).Value = 99.45858346467946
$allTrades.Cells.Item(221, 13).Value = 0
$allTrades.Cells.Item(221, 14).Value = 0
$allTrades.Cells.Item(221, 15).Value = 0
$allTrades.Cells.Item(221, 16).Value = 0.6
$allTrades.Cells.Item(221, 17).Value = "Normal spread capture: 225 bps"

# ---------------------------------------------------------------------------
# "momentum" strategy sheet (row 51 = Trade #190)
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(51, 7).Value = 0.01           # Exit Price
$momentum.Cells.Item(51, 8).Value = "CLOSED"       # Status
$momentum.Cells.Item(51, 11).Value = 99.22         # Capital After
$momentum.Cells.Item(51, 16).Value = "early_exit"  # Exit Reason
$momentum.Cells.Item(51, 17).Value = 0.18          # Duration (min)

# ---------------------------------------------------------------------------
# "HighProbConvergence" strategy sheet - append new trade #219 (row 28)
# ---------------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Range("B28").NumberFormat = "@"
$hpc.Cells.Item(28, 1).Value = 219
$hpc.Cells.Item(28, 2).Value = "2026-02-18"
$hpc.Cells.Item(28, 3).Value = "00:44:21"
$hpc.Cells.Item(28, 4).Value = "HighProbConvergence"
$hpc.Cells.Item(28, 5).Value = "UP"
$hpc.Cells.Item(28, 6).Value = 0.01
$hpc.Cells.Item(28, 8).Value = "OPEN"
$hpc.Cells.Item(28, 9).Value = 0
$hpc.Cells.Item(28, 10).Value = 0
$hpc.Cells.Item(28, 11).Value = 100.3223499536821
$hpc.Cells.Item(28, 12).Value = 0
$hpc.Cells.Item(28, 13).Value = 0
$hpc.Cells.Item(28, 14).Value = 0.95
$hpc.Cells.Item(28, 15).Value = "Mean reversion UP: price 10.82% below mean (z=-3.00)"
$hpc.Cells.Item(28, 17).Value = 0

# ---------------------------------------------------------------------------
# "MarketMaking" strategy sheet - append new trade #220 (row 103)
# ---------------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Range("B103").NumberFormat = "@"
$mm.Cells.Item(103, 1).Value = 220
$mm.Cells.Item(103, 2).Value = "2026-02-18"
$mm.Cells.Item(103, 3).Value = "00:44:22"
$mm.Cells.Item(103, 4).Value = "MarketMaking"
$mm.Cells.Item(103, 5).Value = "UP"
$mm.Cells.Item(103, 6).Value = 0.01
$mm.Cells.Item(103, 8).Value = "OPEN"
$mm.Cells.Item(103, 9).Value = 0
$mm.Cells.Item(103, 10).Value = 0
$mm.Cells.Item(103, 11).Value = 99.45858346467946
$mm.Cells.Item(103, 12).Value = 0
$mm.Cells.Item(103, 13).Value = 0
$mm.Cells.Item(103, 14).Value = 0.6
$mm.Cells.Item(103, 15).Value = "Normal spread capture: 225 bps"
$mm.Cells.Item(103, 17).Value = 0
